$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '40.007.18'
$ws.Range('E2').Value = '  +0.37%  '
Set-TextValue 'D3' '2.237.77'
$ws.Range('E3').Value = '  -4.13%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '294.04'
$ws.Range('E5').Value = '  -4.93%  '
Set-TextValue 'D6' '86.33'
$ws.Range('E6').Value = '  +3.05%  '
Set-TextValue 'D7' '0.515'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.65%  '
Set-TextValue 'D10' '0.0796'
$ws.Range('E10').Value = '  -0.60%  '
Set-TextValue 'D11' '30.54'
$ws.Range('E11').Value = '  +2.83%  '
Set-TextValue 'D12' '47.33'
$ws.Range('E12').Value = '  -9.75%  '
$ws.Range('E13').Value = '  -2.20%  '
Set-TextValue 'D14' '6.39'
$ws.Range('E14').Value = '  +0.21%  '
Set-TextValue 'D15' '2.588.85'
$ws.Range('E15').Value = '  -3.91%  '
Set-TextValue 'D16' '14.23'
$ws.Range('E16').Value = '  -3.06%  '
Set-TextValue 'D17' '2.238.89'
$ws.Range('E17').Value = '  -4.90%  '
Set-TextValue 'D18' '0.726'
$ws.Range('E18').Value = '  -3.42%  '
Set-TextValue 'D19' '39.942.91'
$ws.Range('E19').Value = '  +0.37%  '
Set-TextValue 'D20' '0.0₃0893'
$ws.Range('E20').Value = '  -0.30%  '
Set-TextValue 'D21' '5.80'
$ws.Range('E21').Value = '  -3.80%  '
Set-TextValue 'D22' '10.70'
$ws.Range('E22').Value = '  +2.09%  '
Set-TextValue 'D23' '65.57'
$ws.Range('E23').Value = '  -3.46%  '
Set-TextValue 'D24' '235.14'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  -3.86%  '
Set-TextValue 'D27' '1.84'
$ws.Range('E27').Value = '  +1.81%  '
Set-TextValue 'D28' '23.06'
$ws.Range('E28').Value = '  -1.22%  '
Set-TextValue 'D29' '2.20'
$ws.Range('E29').Value = '  +0.29%  '
Set-TextValue 'D30' '9.25'
$ws.Range('E30').Value = '  +0.42%  '
Set-TextValue 'D31' '33.52'
$ws.Range('E31').Value = '  -0.94%  '
Set-TextValue 'D32' '155.43'
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('E33').Value = '  -0.14%  '
Set-TextValue 'D34' '4.86'
$ws.Range('E34').Value = '  -3.87%  '
Set-TextValue 'D35' '0.0711'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -4.03%  '
Set-TextValue 'D37' '16.54'
$ws.Range('E37').Value = '  +7.17%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  -2.20%  '
Set-TextValue 'D41' '1.67'
$ws.Range('E41').Value = '  -1.50%  '
Set-TextValue 'D42' '3.79'
$ws.Range('E42').Value = '  +0.76%  '
Set-TextValue 'D43' '1.955.71'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('E45').Value = '  +3.22%  '
Set-TextValue 'D46' '9.54'
$ws.Range('E46').Value = '  +1.57%  '
Set-TextValue 'D47' '16.30'
$ws.Range('E47').Value = '  -5.94%  '
Set-TextValue 'D48' '2.62'
$ws.Range('E48').Value = '  -1.49%  '
Set-TextValue 'D49' '2.458.19'
$ws.Range('E49').Value = '  -3.93%  '
Set-TextValue 'D50' '71.02'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('E51').Value = '  +8.59%  '
